$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E are treated as text so values like "56.522.29" or "0.999"
# are not reinterpreted as numbers/dates by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '56.522.29'
$ws.Range("E2").Value = '  +4.70%  '

$ws.Range("D3").Value = '2.995.34'
$ws.Range("E3").Value = '  +5.56%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '508.19'
$ws.Range("E5").Value = '  +10.64%  '

$ws.Range("D6").Value = '137.29'
$ws.Range("E6").Value = '  +11.44%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("E8").Value = '  +7.84%  '

$ws.Range("D9").Value = '7.55'
$ws.Range("E9").Value = '  +15.19%  '

$ws.Range("E10").Value = '  +13.92%  '

$ws.Range("D11").Value = '0.350'
$ws.Range("E11").Value = '  +6.83%  '

$ws.Range("E12").Value = '  +5.88%  '

$ws.Range("D13").Value = '3.507.85'
$ws.Range("E13").Value = '  +5.46%  '

$ws.Range("D14").Value = '25.46'
$ws.Range("E14").Value = '  +10.98%  '

$ws.Range("E15").Value = '  +17.01%  '

$ws.Range("D16").Value = '56.551.56'
$ws.Range("E16").Value = '  +4.70%  '

$ws.Range("D17").Value = '2.999.04'
$ws.Range("E17").Value = '  +5.14%  '

$ws.Range("D18").Value = '5.81'
$ws.Range("E18").Value = '  +9.48%  '

$ws.Range("D19").Value = '12.37'
$ws.Range("E19").Value = '  +9.88%  '

$ws.Range("D20").Value = '7.79'
$ws.Range("E20").Value = '  +12.24%  '

$ws.Range("D21").Value = '326.76'
$ws.Range("E21").Value = '  +11.33%  '

$ws.Range("D22").Value = '0.999'
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("E23").Value = '  +9.84%  '

$ws.Range("D24").Value = '62.37'
$ws.Range("E24").Value = '  +7.82%  '

$ws.Range("D25").Value = '0.170'
$ws.Range("E25").Value = '  +13.96%  '

$ws.Range("E26").Value = '  -0.08%  '

$ws.Range("D27").Value = '0.0₃0909'
$ws.Range("E27").Value = '  +15.81%  '

$ws.Range("D28").Value = '6.55'
$ws.Range("E28").Value = '  +8.68%  '

$ws.Range("D29").Value = '7.06'
$ws.Range("E29").Value = '  +15.56%  '

$ws.Range("D30").Value = '1.26'
$ws.Range("E30").Value = '  +16.96%  '

$ws.Range("E31").Value = '  +13.03%  '

$ws.Range("D32").Value = '20.60'
$ws.Range("E32").Value = '  +11.37%  '

$ws.Range("D33").Value = '155.55'
$ws.Range("E33").Value = '  +13.75%  '

$ws.Range("E34").Value = '  +9.27%  '

$ws.Range("D35").Value = '5.61'
$ws.Range("E35").Value = '  +5.35%  '

$ws.Range("E36").Value = '  +4.93%  '

$ws.Range("E37").Value = '  +10.98%  '

$ws.Range("D38").Value = '24.03'
$ws.Range("E38").Value = '  +6.29%  '

$ws.Range("D39").Value = '3.028.91'
$ws.Range("E39").Value = '  +5.80%  '

$ws.Range("D40").Value = '36.87'
$ws.Range("E40").Value = '  +5.99%  '

$ws.Range("E41").Value = '  +0.00%  '

$ws.Range("D42").Value = '0.646'
$ws.Range("E42").Value = '  +8.07%  '

$ws.Range("D43").Value = '2.258.60'
$ws.Range("E43").Value = '  +11.86%  '

$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").Value = '1.41'
$ws.Range("E44").Value = '  +9.28%  '

$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '0.992'
$ws.Range("E45").Value = '  +7.85%  '

$ws.Range("E46").Value = '  +6.65%  '

$ws.Range("D47").Value = '1.99'
$ws.Range("E47").Value = '  +26.42%  '

$ws.Range("E48").Value = '  +11.60%  '

$ws.Range("E49").Value = '  +8.57%  '

$ws.Range("D50").Value = '19.04'
$ws.Range("E50").Value = '  +9.34%  '

$ws.Range("E51").Value = '  +11.21%  '
